# Add resource(), edit resource(), delete resource() added
#
# The "resource" rows in the sheet shift up one position (a row that used to
# sit further down becomes row 2, etc.), so the "title" cells (column C) end
# up pointing at shared strings that used to carry a trailing space — trim
# them — and the "addTextToTitle" cells (column F) get a "-" prefix on the
# previous "відредагований" text plus a text number format (matching the
# format already used on column E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "Новий ресурс вверху"
$ws.Range("C3").Value = "Новий ресурс в меню"

$ws.Range("F2").Value = "-відредагований"
$ws.Range("F2").NumberFormat = "@"

$ws.Range("F3").Value = "-відредагований"
$ws.Range("F3").NumberFormat = "@"
